$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 195, shifting existing rows 195:203 down to 196:204.
$ws.Rows.Item(195).Insert()

# Populate the new row 195. Columns A,B,C,E,F,G,H,I,J,N,O,Q,R mirror the row that used
# to be at 195 (now at 196); D,K,L,M,P carry the new weekly observation.
$ws.Cells.Item(195, 1).Value = 2
$ws.Cells.Item(195, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(195, 3).Value = "Coquimbo"
$ws.Cells.Item(195, 4).Value = 44706
$ws.Cells.Item(195, 5).Value = 4
$ws.Cells.Item(195, 6).Value = 100112031
$ws.Cells.Item(195, 7).Value = "Poroto verde"
$ws.Cells.Item(195, 8).Value = "Magnum"
$ws.Cells.Item(195, 9).Value = "Primera"
$ws.Cells.Item(195, 10).Value = 400
$ws.Cells.Item(195, 11).Value = 20000
$ws.Cells.Item(195, 12).Value = 23000
$ws.Cells.Item(195, 13).Value = 21500
$ws.Cells.Item(195, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(195, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(195, 16).Value = 860
$ws.Cells.Item(195, 17).Value = 25
$ws.Cells.Item(195, 18).Value = "Hortaliza"

# Match the date-number-format style used by the other date cells in column D.
$ws.Cells.Item(195, 4).NumberFormat = $ws.Cells.Item(196, 4).NumberFormat
